# translations reworked for dino integration
# Rename the German / Italian label columns from locale-tagged names
# ("label::German (de)", "label::Italian (it)") to their ISO3 codes
# ("label::DEU", "label::ITA") on both the "survey" and "choices" sheets.

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item(1)   # sheet1.xml
$choices = $wb.Worksheets.Item(2)   # sheet2.xml

# "survey" sheet: D1 = German label, E1 = Italian label
$survey.Range("D1").Value = "label::DEU"
$survey.Range("E1").Value = "label::ITA"

# "choices" sheet: D1 = Italian label, E1 = German label
$choices.Range("D1").Value = "label::ITA"
$choices.Range("E1").Value = "label::DEU"

# Restore the active-cell selections to E1 on both sheets.
$survey.Activate()
$survey.Range("E1").Select()

$choices.Activate()
$choices.Range("E1").Select()
